$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: replace the SONOVISION ELECTRONICS invoice with the WESTSIDE one.
# D2:H2 must stay plain text (same as the original inlineStr cells), so mark
# them as Text before writing and restore the default "Normal" style
# afterwards so no numFmt/style index leaks onto the cells themselves.
$dataCells = $ws.Range("D2:H2")
$dataCells.NumberFormat = "@"

$ws.Range("A2").Value = "WESTSIDE`nSjr Zion, Survey"
$ws.Range("B2").Value = "29AAACL1838J1ZC"
$ws.Range("C2").Value = "W089 100169940"
$ws.Range("D2").Value = "2024-09-28"
$ws.Range("E2").Value = "4045.01"
$ws.Range("F2").Value = "173.91"
$ws.Range("G2").Value = "173.91"
$ws.Range("H2").Value = "173.91"
$ws.Range("J2").Value = "0`n62052000`n62052000`n62046200`n48194000`n33072000`n39264099"

$dataCells.Style = "Normal"

# Drop the second invoice row (LAKSHMI AGENCIES) entirely; this shifts rows
# up and shrinks the used range back down to A1:J2.
$ws.Rows.Item(3).Delete()
